$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Apollo 13"
$ws.Range("B4").Value = "Start of Construction/build"

$ws.Range("C4").Value = 44586
$ws.Range("D4").Value = 44220
$ws.Range("E4").Value = 44220

$ws.Range("C4:F4").NumberFormat = "dd/mm/yy"
